$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text values keyed by their target shared-string index (2-32): the scraped
# product subtitle (B column, indices 2-18) and price (C column, indices 19-32) HTML snippets.
$strings = @{}
$strings[2] = @'
<p class="product-card__subtitle">
        50% Wool 50% Acrylic, 180m (197yds)/100g (3.5oz), Aran
      </p>
'@
$strings[3] = @'
<p class="product-card__subtitle">
        100% Wool, 88m (96yds)/50g (1.8oz), Aran
      </p>
'@
$strings[4] = @'
<p class="product-card__subtitle">
        100% Wool, 880m (962yds)/500g (17.6oz), Aran
      </p>
'@
$strings[5] = @'
<p class="product-card__subtitle">
        100% Wool, 200m (219yds)/100g (3.5oz), Aran
      </p>
'@
$strings[6] = @'
<p class="product-card__subtitle">
        100% Wool, 166m (182yds)/100g (3.5oz), Aran
      </p>
'@
$strings[7] = @'
<p class="product-card__subtitle">
        100% Wool, 90m (98yds)/50g (1.8oz), Aran
      </p>
'@
$strings[8] = @'
<p class="product-card__subtitle">
        100% Wool, 900m (984yds)/500g (17.6oz), Aran
      </p>
'@
$strings[9] = @'
<p class="product-card__subtitle">
        100% Wool, 440m (481yds)/250g (8.8oz), Aran
      </p>
'@
$strings[10] = @'
<p class="product-card__subtitle">
        50% Wool 25% Alpaca 25% Viscose, 87m (95yds)/50g (1.8oz), Aran
      </p>
'@
$strings[11] = @'
<p class="product-card__subtitle">
        100% Wool, 100m (109yds)/50g (1.8oz), Aran
      </p>
'@
$strings[12] = @'
<p class="product-card__subtitle">
        100% Wool, 180m (197yds)/100g (3.5oz), Aran
      </p>
'@
$strings[13] = @'
<p class="product-card__subtitle">
        50% Wool 50% Acrylic, 500g (17.6oz)
      </p>
'@
$strings[14] = @'
<p class="product-card__subtitle">
        35% Wool 20% Acrylic 20% Polyamide 25% Viscose, 300m (328yds)/100g (3.5oz), Aran
      </p>
'@
$strings[15] = @'
<p class="product-card__subtitle">
        100% Wool, 137m (150yds)/100g (3.5oz), Aran
      </p>
'@
$strings[16] = @'
<p class="product-card__subtitle">
        96% Wool 4% Viscose, 166m (182yds)/100g (3.5oz), Aran
      </p>
'@
$strings[17] = @'
<p class="product-card__subtitle">
        80% Wool 20% Polyamide, 199m (218yds)/100g (3.5oz), Aran
      </p>
'@
$strings[18] = @'
<p class="product-card__subtitle">
</p>
'@
$strings[19] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £7.00
    </span></span>
'@
$strings[20] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £8.49
    </span></span>
'@
$strings[21] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £9.99
    </span></span>
'@
$strings[22] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £7.49
    </span></span>
'@
$strings[23] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £9.49
    </span></span>
'@
$strings[24] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £6.95
    </span></span>
'@
$strings[25] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £4.99
    </span></span>
'@
$strings[26] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £3.49
    </span></span>
'@
$strings[27] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £8.70
    </span></span>
'@
$strings[28] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £6.90
    </span></span>
'@
$strings[29] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £4.49
    </span></span>
'@
$strings[30] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £10.99
    </span></span>
'@
$strings[31] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £12.99
    </span></span>
'@
$strings[32] = @'
<span class="lc-price__regular" data-v-27ab4212=""><span data-v-27ab4212="">
      £6.25
    </span></span>
'@

# Row data: (row, A value, B shared-string index, C shared-string index)
$rows = @(
    @(2, 0, 2, 19),
    @(3, 1, 3, 20),
    @(4, 2, 4, 21),
    @(5, 3, 5, 22),
    @(6, 4, 6, 23),
    @(7, 5, 7, 24),
    @(8, 6, 8, 20),
    @(9, 7, 9, 25),
    @(10, 8, 6, 26),
    @(11, 9, 10, 27),
    @(12, 10, 6, 24),
    @(13, 11, 11, 28),
    @(14, 12, 12, 29),
    @(15, 13, 13, 21),
    @(16, 14, 14, 30),
    @(17, 15, 15, 31),
    @(18, 16, 16, 23),
    @(19, 17, 17, 32),
    @(20, 18, 18, 21)
)

# Write column B (product subtitle) for every row first, then column C (price) for
# every row, then column A — this matches the shared-string insertion order of the target
# workbook (all subtitles, then all prices).
foreach ($entry in $rows) {
    $r = $entry[0]
    $bIdx = $entry[2]
    $ws.Cells.Item($r, 2).Value = $strings[$bIdx]
}

foreach ($entry in $rows) {
    $r = $entry[0]
    $cIdx = $entry[3]
    $ws.Cells.Item($r, 3).Value = $strings[$cIdx]
}

foreach ($entry in $rows) {
    $r = $entry[0]
    $aVal = $entry[1]

    # Column A gets the bold/bordered/centered header style (copied from B1) plus its numeric value
    $ws.Range("B1").Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $aVal
}
$excel.CutCopyMode = $false

# The multi-line HTML snippets bump each rows auto height; AutoFit recomputes it back
# to the sheet default (15) since WrapText is off, matching the target (no ht override).
foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Rows.Item($r).EntireRow.AutoFit()
}

Write-Output "done"
